$wb = $excel.ActiveWorkbook

# 1. Add a new blank worksheet ("Sheet1") positioned right after "mainTimeline"
$newSheet1 = $wb.Worksheets.Add()
$newSheet1.Name = "Sheet1"
$after = $wb.Worksheets.Item("mainTimeline")
$newSheet1.Move($null, $after)

# 2. Duplicate "usefulInfo" sheet (copy after itself), then rename copy to "Sheet2"
#    and move it to sit right after "studyDesignOE"
$usefulInfo = $wb.Worksheets.Item("usefulInfo")
$usefulInfo.Copy($null, $usefulInfo)
$newSheet2 = $wb.ActiveSheet
$newSheet2.Name = "Sheet2"
$oe = $wb.Worksheets.Item("studyDesignOE")
$newSheet2.Move($null, $oe)

# Re-fetch the sheet by its (now-stable) name, since Move() can invalidate
# the previously-held object reference, then make it the active/selected tab
$newSheet2 = $wb.Worksheets.Item("Sheet2")
$newSheet2.Select()
